$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for each changed cell so numeric-looking strings
# (e.g. "582.63") are not auto-converted to real numbers, matching the
# original inlineStr/text representation used throughout column D and E.
function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue $ws 'D2' '67.111.63'
Set-TextValue $ws 'E2' '  -0.05%  '
Set-TextValue $ws 'D3' '2.467.74'
Set-TextValue $ws 'E3' '  -0.14%  '
Set-TextValue $ws 'E4' '  +0.02%  '
Set-TextValue $ws 'D5' '582.63'
Set-TextValue $ws 'E5' '  -0.10%  '
Set-TextValue $ws 'D6' '174.14'
Set-TextValue $ws 'E6' '  +3.03%  '
Set-TextValue $ws 'E7' '  +0.07%  '
Set-TextValue $ws 'E8' '  -0.43%  '
Set-TextValue $ws 'E9' '  +1.72%  '
Set-TextValue $ws 'E10' '  +0.48%  '
Set-TextValue $ws 'E11' '  -0.68%  '
Set-TextValue $ws 'E12' '  +0.44%  '
Set-TextValue $ws 'D13' '2.928.67'
Set-TextValue $ws 'E13' '  +2.17%  '
Set-TextValue $ws 'D14' '25.35'
Set-TextValue $ws 'E14' '  -0.98%  '
Set-TextValue $ws 'D15' '66.998.03'
Set-TextValue $ws 'E15' '  -0.08%  '
Set-TextValue $ws 'E16' '  -0.31%  '
Set-TextValue $ws 'D17' '2.457.44'
Set-TextValue $ws 'E17' '  +0.87%  '
Set-TextValue $ws 'E18' '  -2.52%  '
Set-TextValue $ws 'D19' '7.46'
Set-TextValue $ws 'E19' '  -2.00%  '
Set-TextValue $ws 'D20' '348.62'
Set-TextValue $ws 'E21' '  -1.49%  '
Set-TextValue $ws 'E22' '  +0.00%  '
Set-TextValue $ws 'D23' '69.26'
Set-TextValue $ws 'E23' '  +0.28%  '
Set-TextValue $ws 'E24' '  -1.22%  '
Set-TextValue $ws 'E25' '  +0.34%  '
Set-TextValue $ws 'E26' '  -1.86%  '
Set-TextValue $ws 'D27' '2.594.65'
Set-TextValue $ws 'E27' '  +0.21%  '
Set-TextValue $ws 'D28' '0.997'
Set-TextValue $ws 'E28' '  -0.02%  '
Set-TextValue $ws 'D29' '0.0₃0898'
Set-TextValue $ws 'E29' '  -1.06%  '
Set-TextValue $ws 'D30' '499.38'
Set-TextValue $ws 'E30' '  -3.35%  '
Set-TextValue $ws 'D31' '7.71'
Set-TextValue $ws 'E31' '  -0.55%  '
Set-TextValue $ws 'E32' '  -1.19%  '
Set-TextValue $ws 'E33' '  -1.69%  '
Set-TextValue $ws 'E34' '  +0.00%  '
Set-TextValue $ws 'D35' '0.119'
Set-TextValue $ws 'E35' '  -0.25%  '
Set-TextValue $ws 'D36' '161.60'
Set-TextValue $ws 'E36' '  +2.47%  '
Set-TextValue $ws 'D38' '18.11'
Set-TextValue $ws 'E38' '  -1.55%  '
Set-TextValue $ws 'E39' '  -1.99%  '
Set-TextValue $ws 'E40' '  -0.01%  '
Set-TextValue $ws 'D41' '1.68'
Set-TextValue $ws 'E41' '  +1.38%  '
Set-TextValue $ws 'D43' '4.81'
Set-TextValue $ws 'E43' '  +0.39%  '
Set-TextValue $ws 'E44' '  +0.20%  '
Set-TextValue $ws 'D45' '141.97'
Set-TextValue $ws 'E45' '  +0.58%  '
Set-TextValue $ws 'D46' '3.46'
Set-TextValue $ws 'E46' '  +0.06%  '
Set-TextValue $ws 'D47' '0.510'
Set-TextValue $ws 'E47' '  -1.14%  '
Set-TextValue $ws 'D48' '0.0₆0254'
Set-TextValue $ws 'E48' '  -0.37%  '
Set-TextValue $ws 'E49' '  +0.34%  '
Set-TextValue $ws 'E50' '  -1.67%  '
Set-TextValue $ws 'D51' '0.581'
Set-TextValue $ws 'E51' '  +0.03%  '
